# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet between the existing "2021-Q4" sheet
#   and the "总计" (totals) sheet, populated with the quarter's fund-holder
#   detail (same layout/headers as "2021-Q4").
# - Prepends a "2022-Q1" row to the "总计" sheet's summary table, pushing the
#   existing "2021-Q4" summary row down by one.

$wb = $excel.ActiveWorkbook

$detailSheet = $wb.Worksheets.Item(1)   # "2021-Q4" - template layout to copy
$totalSheetBefore = $wb.Worksheets.Item(2)   # "总计" (pre-insert position)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right before "总计" (so tab order is
#    2021-Q4, 2022-Q1, 总计) and clone the header/row formatting from the
#    "2021-Q4" sheet.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"
# Match the outline summary-direction flags the other sheets in this
# workbook carry (summary rows below / summary columns to the right).
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# NOTE: "Worksheets.Item(N)" resolves by live position, so the reference
# captured above as $totalSheetBefore now actually points at $newSheet
# (whatever sheet currently sits at position 2). Re-resolve "总计" by name
# now that the new sheet has been inserted ahead of it.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row + sample data row (columns B:H) - values, then formats.
$detailSheet.Range("B1:H2").Copy()
$newSheet.Range("B1").PasteSpecial(-4104)
$detailSheet.Range("B1:H2").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

# Column A (the bold/bordered running-index column) only has a formatted
# cell on the data row, not the header row - copy its format separately.
$detailSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Overwrite the copied sample row with the real 2022-Q1 fund data. The
# numeric-looking values are stored as text (matching the source sheet),
# only the rank column (H) is a genuine number.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'168207"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "中融创业板两年定期开放混合"
$newSheet.Range("D2").Value = "'3.85"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'88.46"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'3.85"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.1482"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 5

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet's data,
#    shifting "2021-Q4" down to row 3 (preserving its values/format).
# ---------------------------------------------------------------------
$oldDate  = $totalSheet.Range("B2").Value()
$oldCount = $totalSheet.Range("C2").Value()
$oldValue = $totalSheet.Range("D2").Value()

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = $oldDate
$totalSheet.Range("C3").Value = $oldCount
$totalSheet.Range("D3").Value = $oldValue

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.15
